# Update gh-pages to output generated at 456a3b4
# Refresh "want-to-go" counts (column F) and a couple of price/status
# cells (column G) across the four sheets of the workbook.

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (exhibitions) ----
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 4586
$ws.Range("F5").Value  = 3691
$ws.Range("F6").Value  = 1074
$ws.Range("F10").Value = 374
$ws.Range("F11").Value = 2571
$ws.Range("F12").Value = 1291
$ws.Range("F15").Value = 280
$ws.Range("F16").Value = 23
$ws.Range("F17").Value = 564
$ws.Range("F18").Value = 265
$ws.Range("F20").Value = 10655
$ws.Range("F21").Value = 6152
$ws.Range("F29").Value = 26
$ws.Range("F30").Value = 190
$ws.Range("F31").Value = 866
$ws.Range("F32").Value = 3573
$ws.Range("F34").Value = 972
$ws.Range("F36").Value = 135
$ws.Range("F39").Value = 256
$ws.Range("F40").Value = 4874
$ws.Range("F42").Value = 1157
$ws.Range("F44").Value = 206
$ws.Range("F45").Value = 120
$ws.Range("F46").Value = 496

# ---- Sheet "演出" (performances) ----
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value  = "不可售"
$ws.Range("F7").Value  = 16
$ws.Range("F11").Value = 89
$ws.Range("F15").Value = 3610
$ws.Range("G15").Value = 398

# ---- Sheet "本地生活" (local life) ----
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 8858
$ws.Range("F4").Value = 1675

# ---- Sheet "全部类型" (all types) ----
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 8858
$ws.Range("F4").Value  = 1675
$ws.Range("F5").Value  = 4586
$ws.Range("F7").Value  = 3691
$ws.Range("F8").Value  = 1074
$ws.Range("F11").Value = 374
$ws.Range("F12").Value = 2571
$ws.Range("F13").Value = 16
$ws.Range("F16").Value = 1291
$ws.Range("F17").Value = 89
$ws.Range("F19").Value = 280
$ws.Range("F20").Value = 23
$ws.Range("F21").Value = 564
$ws.Range("F22").Value = 265
$ws.Range("F23").Value = 10655
$ws.Range("F24").Value = 3610
$ws.Range("G24").Value = 398
$ws.Range("F32").Value = 26
$ws.Range("F33").Value = 866
$ws.Range("F34").Value = 3573
$ws.Range("F36").Value = 972
$ws.Range("F37").Value = 135
$ws.Range("F41").Value = 256
$ws.Range("F42").Value = 4874
$ws.Range("F44").Value = 1157
$ws.Range("F46").Value = 120
$ws.Range("F47").Value = 496
